# Commit: "Push choices sheet display.text into display.title.text"
# (Also push display.image into display.title.image -- not applicable to
#  this workbook, which has no display.image column.)
#
# The "choices" sheet's header row currently labels column C as
# "display.text". Rename that header to "display.title.text" so the
# choice-list display text is driven from the title field instead.

$wb = $excel.ActiveWorkbook

$choices = $wb.Worksheets.Item("choices")
$choices.Range("C1").Value = "display.title.text"

# Reflect the real author's interactive session: they ended up on the
# choices tab with C8 selected.
$choices.Activate()
$choices.Range("C8").Select()
